$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7; existing rows 7-52 shift down to 8-53.
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new "Té" (tea) product record.
$ws.Cells.Item(7, 1).Value = 7790150211953
$ws.Cells.Item(7, 2).Value = "Té"
$ws.Cells.Item(7, 3).Value = "negro clásico"
$ws.Cells.Item(7, 4).Value = "en saquitos"
$ws.Cells.Item(7, 5).Value = "La Virginia"
$ws.Cells.Item(7, 6).Value = 100
$ws.Cells.Item(7, 7).Value = "und."
$ws.Cells.Item(7, 8).Value = "caja"
$ws.Cells.Item(7, 9).Value = "Tés"
$ws.Cells.Item(7, 10).Value = "Argentina"
$ws.Cells.Item(7, 11).Value = 6
$ws.Cells.Item(7, 12).Value = $false
$ws.Cells.Item(7, 13).Value = $true
$ws.Cells.Item(7, 14).Value = "C:\VentaSoft\Imágenes de artículos\7790150211953.png"
$ws.Cells.Item(7, 16).Value = $false

# Column O (ImagenExactaDelArticulo) needs the same boolean TRUE value and the
# exact cell style used by the surrounding rows (copy it straight from the row
# below, which still carries the original "Normal_Artículos" style).
$ws.Cells.Item(8, 15).Copy($ws.Cells.Item(7, 15))
